# Insert a new weekly price record for "Haba" (Macroferia Regional de Talca)
# as row 85, pushing the existing rows 85-140 down to 86-141.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(85).Insert()

$ws.Range("A85").Value = 5
$ws.Range("B85").Value = "Macroferia Regional de Talca"
$ws.Range("C85").Value = "Maule"
$ws.Range("D85").Value = 45233
$ws.Range("E85").Value = 7
$ws.Range("F85").Value = 100112026
$ws.Range("G85").Value = "Haba"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 500
$ws.Range("K85").Value = 9000
$ws.Range("L85").Value = 9000
$ws.Range("M85").Value = 9000
$ws.Range("N85").Value = "$/saco 25 kilos"
$ws.Range("O85").Value = "Región del Maule"
$ws.Range("P85").Value = 360
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"
